# Guardians Orders character table: add a 4-column "Skill Set" block
# (SkillQ / SkillW / SkillE / SkillR) right after the character-name
# column, pushing the existing Default*/Skeletal/Anim columns to the
# right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new (empty) columns before the old column B, i.e. columns B:E.
# This shifts the former B..H block to F..L.
$ws.Range("B1:E1").EntireColumn.Insert()

# The insert cloned column A's formatting/style down every row of B:E.
# Only row 1 should carry the new headers; rows 2:5 must go back to being
# completely empty (no cell entries at all), matching the rest of the
# (data-only, header-less) new columns.
$ws.Range("B2:E5").Clear()

# New header row values for the inserted columns.
$ws.Range("B1").Value = "SkillQ"
$ws.Range("C1").Value = "SkillW"
$ws.Range("D1").Value = "SkillE"
$ws.Range("E1").Value = "SkillR"

# Column widths: A:E become a uniform narrower width, while F (old B) and
# G (old C) are slightly resized too.
$ws.Range("A1:E1").EntireColumn.ColumnWidth = 15.571428571428571
$ws.Range("F1").EntireColumn.ColumnWidth = 23.42857142857143
$ws.Range("G1").EntireColumn.ColumnWidth = 19.857142857142854

# Selection moves to the first new header cell.
$ws.Range("E1").Select() | Out-Null
